# Matriz de Trazabilidad - actualizacion de estado y entregables
# (5to trimestre: se agregan documentos / codigo fuente, requisitos pasan a "En curso")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Filas 4 a 13 contienen los requerimientos (RF01, RF03-RF11)
# Columna F = "Estado actual"      -> pasa de "Incompleto" a "En curso"
# Columna J = "Entregables"        -> pasa de "N/A" a "Codigo fuente"
for ($row = 4; $row -le 13; $row++) {
    $ws.Range("F$row").Value = "En curso"
    $ws.Range("J$row").Value = "Codigo fuente"
}

# El contenido de la columna J es ahora mas largo ("Codigo fuente"), se reajusta el ancho
$ws.Columns("J:J").AutoFit() | Out-Null

# Se deja seleccionado el rango de entregables recien actualizado
$ws.Range("J4:J13").Select() | Out-Null
